# Generate Report for Handoff
# Adds two new localization entries (805f917c... and e924d709...) to the
# Overview / zh-cn / de-de sheets, reorders the first two entries
# (06159bb4... now before 2bab85a4...), and pushes the existing
# ".localization-config" row down accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Drop existing hyperlinks so we can rebuild them in the new row order
# without leaving stale entries behind (Hyperlinks.Delete on any Range
# clears the whole sheet's collection in this engine).
$ov.Hyperlinks.Delete()

$ov.Range("A2").Value = "06159bb4-5dff-4fa3-b03a-181ceb630147.md"
$ov.Range("B2").Value = "In Translation"
$ov.Range("C2").Value = "In Translation"

$ov.Range("A3").Value = "2bab85a4-5ea1-4745-806b-4952705b6eb0.md"
$ov.Range("B3").Value = "In Translation"
$ov.Range("C3").Value = "In Translation"

$ov.Range("A4").Value = "805f917c-65e4-40c3-bbee-9b213eb88d4c.md"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"

$ov.Range("A5").Value = "e924d709-8e56-4d53-bdb9-cc448585afe6.md"
$ov.Range("B5").Value = "Ready for handoff"
$ov.Range("C5").Value = "Ready for handoff"

$ov.Range("A6").Value = ".localization-config"
$ov.Range("B6").Value = "Not to be localized"
$ov.Range("C6").Value = "Not to be localized"

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0de453485eb110f381433ba6278e850be0f4e3b6/e2e/06159bb4-5dff-4fa3-b03a-181ceb630147.md", "", "", "06159bb4-5dff-4fa3-b03a-181ceb630147.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0de453485eb110f381433ba6278e850be0f4e3b6/e2e/2bab85a4-5ea1-4745-806b-4952705b6eb0.md", "", "", "2bab85a4-5ea1-4745-806b-4952705b6eb0.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0de453485eb110f381433ba6278e850be0f4e3b6/e2e/805f917c-65e4-40c3-bbee-9b213eb88d4c.md", "", "", "805f917c-65e4-40c3-bbee-9b213eb88d4c.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0de453485eb110f381433ba6278e850be0f4e3b6/e2e/e924d709-8e56-4d53-bdb9-cc448585afe6.md", "", "", "e924d709-8e56-4d53-bdb9-cc448585afe6.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/0de453485eb110f381433ba6278e850be0f4e3b6/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Hyperlinks.Delete()

$zh.Range("A2").Value = "06159bb4-5dff-4fa3-b03a-181ceb630147.md"
$zh.Range("B2").Value = "In Translation"
$zh.Range("C2").Value = "06159bb4-5dff-4fa3-b03a-181ceb630147.4c1195f4a01bfca3a06184e4124bb9795a800f94.zh-cn.xlf"
$zh.Range("D2").Value = "2016-03-03 12:12:55"
$zh.Range("G2").Value = "0001-01-01 00:00:00"
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = "2bab85a4-5ea1-4745-806b-4952705b6eb0.md"
$zh.Range("B3").Value = "In Translation"
$zh.Range("C3").Value = "2bab85a4-5ea1-4745-806b-4952705b6eb0.85d85ef246863b7cc26f5f350f7e2f4b4b35acb4.zh-cn.xlf"
$zh.Range("D3").Value = "2016-03-03 12:12:55"
$zh.Range("G3").Value = "0001-01-01 00:00:00"
$zh.Range("H3").Value = "Include"

$zh.Range("A4").Value = "805f917c-65e4-40c3-bbee-9b213eb88d4c.md"
$zh.Range("B4").Value = "Ready for handoff"
$zh.Range("C4").Value = "805f917c-65e4-40c3-bbee-9b213eb88d4c.232e74c988645b1f968838d73c84a323c2f3c410.zh-cn.xlf"
$zh.Range("D4").Value = "2016-03-03 12:15:02"
$zh.Range("G4").Value = "0001-01-01 00:00:00"
$zh.Range("H4").Value = "Include"

$zh.Range("A5").Value = "e924d709-8e56-4d53-bdb9-cc448585afe6.md"
$zh.Range("B5").Value = "Ready for handoff"
$zh.Range("C5").Value = "e924d709-8e56-4d53-bdb9-cc448585afe6.515c5dece074ea49ee26839ffbc295103c0ade0f.zh-cn.xlf"
$zh.Range("D5").Value = "2016-03-03 12:15:02"
$zh.Range("G5").Value = "0001-01-01 00:00:00"
$zh.Range("H5").Value = "Include"

$zh.Range("A6").Value = ".localization-config"
$zh.Range("B6").Value = "Not to be localized"
$zh.Range("D6").Value = "0001-01-01 00:00:00"
$zh.Range("G6").Value = "0001-01-01 00:00:00"
$zh.Range("H6").Value = "Ignored"

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0de453485eb110f381433ba6278e850be0f4e3b6/e2e/06159bb4-5dff-4fa3-b03a-181ceb630147.md", "", "", "06159bb4-5dff-4fa3-b03a-181ceb630147.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/675bc1a0c906a83954d1e2841842a8a7e2442bb4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/06159bb4-5dff-4fa3-b03a-181ceb630147.4c1195f4a01bfca3a06184e4124bb9795a800f94.zh-cn.xlf", "", "", "06159bb4-5dff-4fa3-b03a-181ceb630147.4c1195f4a01bfca3a06184e4124bb9795a800f94.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0de453485eb110f381433ba6278e850be0f4e3b6/e2e/2bab85a4-5ea1-4745-806b-4952705b6eb0.md", "", "", "2bab85a4-5ea1-4745-806b-4952705b6eb0.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/675bc1a0c906a83954d1e2841842a8a7e2442bb4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/2bab85a4-5ea1-4745-806b-4952705b6eb0.85d85ef246863b7cc26f5f350f7e2f4b4b35acb4.zh-cn.xlf", "", "", "2bab85a4-5ea1-4745-806b-4952705b6eb0.85d85ef246863b7cc26f5f350f7e2f4b4b35acb4.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0de453485eb110f381433ba6278e850be0f4e3b6/e2e/805f917c-65e4-40c3-bbee-9b213eb88d4c.md", "", "", "805f917c-65e4-40c3-bbee-9b213eb88d4c.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/675bc1a0c906a83954d1e2841842a8a7e2442bb4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/805f917c-65e4-40c3-bbee-9b213eb88d4c.232e74c988645b1f968838d73c84a323c2f3c410.zh-cn.xlf", "", "", "805f917c-65e4-40c3-bbee-9b213eb88d4c.232e74c988645b1f968838d73c84a323c2f3c410.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0de453485eb110f381433ba6278e850be0f4e3b6/e2e/e924d709-8e56-4d53-bdb9-cc448585afe6.md", "", "", "e924d709-8e56-4d53-bdb9-cc448585afe6.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/675bc1a0c906a83954d1e2841842a8a7e2442bb4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/e924d709-8e56-4d53-bdb9-cc448585afe6.515c5dece074ea49ee26839ffbc295103c0ade0f.zh-cn.xlf", "", "", "e924d709-8e56-4d53-bdb9-cc448585afe6.515c5dece074ea49ee26839ffbc295103c0ade0f.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/0de453485eb110f381433ba6278e850be0f4e3b6/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Hyperlinks.Delete()

$de.Range("A2").Value = "06159bb4-5dff-4fa3-b03a-181ceb630147.md"
$de.Range("B2").Value = "In Translation"
$de.Range("C2").Value = "06159bb4-5dff-4fa3-b03a-181ceb630147.4c1195f4a01bfca3a06184e4124bb9795a800f94.de-de.xlf"
$de.Range("D2").Value = "2016-03-03 12:13:10"
$de.Range("G2").Value = "0001-01-01 00:00:00"
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = "2bab85a4-5ea1-4745-806b-4952705b6eb0.md"
$de.Range("B3").Value = "In Translation"
$de.Range("C3").Value = "2bab85a4-5ea1-4745-806b-4952705b6eb0.85d85ef246863b7cc26f5f350f7e2f4b4b35acb4.de-de.xlf"
$de.Range("D3").Value = "2016-03-03 12:13:10"
$de.Range("G3").Value = "0001-01-01 00:00:00"
$de.Range("H3").Value = "Include"

$de.Range("A4").Value = "805f917c-65e4-40c3-bbee-9b213eb88d4c.md"
$de.Range("B4").Value = "Ready for handoff"
$de.Range("C4").Value = "805f917c-65e4-40c3-bbee-9b213eb88d4c.232e74c988645b1f968838d73c84a323c2f3c410.de-de.xlf"
$de.Range("D4").Value = "2016-03-03 12:15:15"
$de.Range("G4").Value = "0001-01-01 00:00:00"
$de.Range("H4").Value = "Include"

$de.Range("A5").Value = "e924d709-8e56-4d53-bdb9-cc448585afe6.md"
$de.Range("B5").Value = "Ready for handoff"
$de.Range("C5").Value = "e924d709-8e56-4d53-bdb9-cc448585afe6.515c5dece074ea49ee26839ffbc295103c0ade0f.de-de.xlf"
$de.Range("D5").Value = "2016-03-03 12:15:15"
$de.Range("G5").Value = "0001-01-01 00:00:00"
$de.Range("H5").Value = "Include"

$de.Range("A6").Value = ".localization-config"
$de.Range("B6").Value = "Not to be localized"
$de.Range("D6").Value = "0001-01-01 00:00:00"
$de.Range("G6").Value = "0001-01-01 00:00:00"
$de.Range("H6").Value = "Ignored"

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0de453485eb110f381433ba6278e850be0f4e3b6/e2e/06159bb4-5dff-4fa3-b03a-181ceb630147.md", "", "", "06159bb4-5dff-4fa3-b03a-181ceb630147.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b66c66fd04380bbcaba53b40dbf3cc08b739c32f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/06159bb4-5dff-4fa3-b03a-181ceb630147.4c1195f4a01bfca3a06184e4124bb9795a800f94.de-de.xlf", "", "", "06159bb4-5dff-4fa3-b03a-181ceb630147.4c1195f4a01bfca3a06184e4124bb9795a800f94.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0de453485eb110f381433ba6278e850be0f4e3b6/e2e/2bab85a4-5ea1-4745-806b-4952705b6eb0.md", "", "", "2bab85a4-5ea1-4745-806b-4952705b6eb0.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b66c66fd04380bbcaba53b40dbf3cc08b739c32f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/2bab85a4-5ea1-4745-806b-4952705b6eb0.85d85ef246863b7cc26f5f350f7e2f4b4b35acb4.de-de.xlf", "", "", "2bab85a4-5ea1-4745-806b-4952705b6eb0.85d85ef246863b7cc26f5f350f7e2f4b4b35acb4.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0de453485eb110f381433ba6278e850be0f4e3b6/e2e/805f917c-65e4-40c3-bbee-9b213eb88d4c.md", "", "", "805f917c-65e4-40c3-bbee-9b213eb88d4c.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b66c66fd04380bbcaba53b40dbf3cc08b739c32f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/805f917c-65e4-40c3-bbee-9b213eb88d4c.232e74c988645b1f968838d73c84a323c2f3c410.de-de.xlf", "", "", "805f917c-65e4-40c3-bbee-9b213eb88d4c.232e74c988645b1f968838d73c84a323c2f3c410.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0de453485eb110f381433ba6278e850be0f4e3b6/e2e/e924d709-8e56-4d53-bdb9-cc448585afe6.md", "", "", "e924d709-8e56-4d53-bdb9-cc448585afe6.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b66c66fd04380bbcaba53b40dbf3cc08b739c32f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/e924d709-8e56-4d53-bdb9-cc448585afe6.515c5dece074ea49ee26839ffbc295103c0ade0f.de-de.xlf", "", "", "e924d709-8e56-4d53-bdb9-cc448585afe6.515c5dece074ea49ee26839ffbc295103c0ade0f.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/0de453485eb110f381433ba6278e850be0f4e3b6/.localization-config", "", "", ".localization-config") | Out-Null
